# Logged Week 15 and simulated Week 16
# Update the "H" (home) row target-depth stats on both the OFF and DEF sheets
# with the latest cumulative season totals.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 124
$wsOff.Range("C2").Value = 80
$wsOff.Range("D2").Value = 41
$wsOff.Range("E2").Value = 20
$wsOff.Range("F2").Value = 1

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 175
$wsDef.Range("C2").Value = 124
$wsDef.Range("D2").Value = 31
$wsDef.Range("E2").Value = 13
